# Applies the edits described by the diff to the active document.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Trim the long justification paragraph down to a single trailing
#    space (removes the "O dono..." / "Para alugar..." sentences).
# ---------------------------------------------------------------------
$old1 = " O dono terá acesso a todas funcionalidades do sistema e deverá controlar os níveis de acesso. Para alugar filmes e jogos, o usuário tem que ser previamente cadastrado no sistema. O cadastro é na hora do pagamento caso ele não esteja cadastrado, necessitando apresentar o nome, cpf e telefone. "
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, " ", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Remove the "_GoBack" bookmark that currently sits right after
#    "o distribuidor".
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 3) Drop "localiza o cadastro do usuário no sistema, " from the
#    rental paragraph.
# ---------------------------------------------------------------------
$old3 = "Ao alugar um exemplar, o funcionário localiza o cadastro do usuário no sistema, verifica a situação do usuário na relação de exemplares"
$new3 = "Ao alugar um exemplar, o funcionário verifica a situação do usuário na relação de exemplares"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Rewrite the exclusions paragraph opening and re-insert the
#    "_GoBack" bookmark right before "e pode ocorrer por mot".
# ---------------------------------------------------------------------
$old4 = "As exclusões de usuários são informadas pelo dono. Já a exclusão de exemplares é decidida pelo funcionário, mas precisa de autorização do dono e pode ocorrer por mot"
$new4 = "A exclusão de exemplares é decidida pelo funcionário, e pode ocorrer por mot"
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

$anchor = "e pode ocorrer por mot"
$r = $d.Content
$r.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pos = $r.Start
$bm = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bm) | Out-Null
